$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$s.Shapes.Item(1).TextFrame.TextRange.Text = "GSoC - 2018 Project Proposal"
$s.Shapes.Item(2).TextFrame.TextRange.Text = "Akshay Anand"
